$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment the "Förändrad" date (column C) from 45650 to 45651 for rows 2-37
for ($r = 2; $r -le 37; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45650) {
        $cell.Value2 = 45651
    }
}

# Swap the "Beteckning" (column A) values between rows 36 and 37
$a36 = $ws.Cells.Item(36, 1).Value2
$a37 = $ws.Cells.Item(37, 1).Value2
$ws.Cells.Item(36, 1).Value2 = $a37
$ws.Cells.Item(37, 1).Value2 = $a36

# Swap the "Area (ha)" (column G) values between rows 36 and 37
$g36 = $ws.Cells.Item(36, 7).Value2
$g37 = $ws.Cells.Item(37, 7).Value2
$ws.Cells.Item(36, 7).Value2 = $g37
$ws.Cells.Item(37, 7).Value2 = $g36
